$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value as TEXT (not auto-converted to a number) while
# keeping the cell's style at the sheet's plain/default style (no "s" index).
# We do this by temporarily forcing a text number-format, writing the value,
# then pasting the (blank) format from a never-touched "donor" cell back on
# top so the final style returns to the default.
# ---------------------------------------------------------------------------
function Set-TextValue($cell, $donor, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $donor.Copy()
    $cell.PasteSpecial(-4122)
}

# ===========================================================================
# 1) Insert the new "2022-Q4" detail sheet right after "总计".
# ===========================================================================
$summary = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $summary)
$q4.Name = "2022-Q4"
$donorQ4 = $q4.Cells.Item(100, 100)

# Header row (style copied from an existing header cell so it gets the bold
# "s=2" look used throughout the workbook).
$hdrDonor = $wb.Worksheets.Item(3).Cells.Item(1, 2)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $q4.Cells.Item(1, $col)
    $hdrDonor.Copy()
    $cell.PasteSpecial(-4122)
    $cell.Value = $headers[$col - 2]
}

# Data rows
$rows = @(
    @("270023", "广发全球精选股票（QDII）", "20.45", "82.63", "6.97", "1.4254", 4),
    @("000906", "广发全球精选股票（QDII）美元现汇", "20.45", "82.63", "6.97", "1.4254", 4),
    @("015203", "汇添富全球移动互联灵活配置混合（QDII）D", "11.52", "92.14", "3.48", "0.4009", 5),
    @("001668", "汇添富全球移动互联灵活配置混合（QDII）A", "11.48", "92.14", "3.48", "0.3995", 5),
    @("015202", "汇添富全球移动互联灵活配置混合（QDII）C", "0.01", "92.14", "3.48", "0.0003", 5)
)

$rowDonor = $wb.Worksheets.Item(3).Cells.Item(2, 1)
$r = 2
foreach ($row in $rows) {
    $rowDonor.Copy()
    $q4.Cells.Item($r, 1).PasteSpecial(-4122)
    $q4.Cells.Item($r, 1).Value = ($r - 2)

    Set-TextValue $q4.Cells.Item($r, 2) $donorQ4 $row[0]
    $q4.Cells.Item($r, 3).Value = $row[1]
    Set-TextValue $q4.Cells.Item($r, 4) $donorQ4 $row[2]
    Set-TextValue $q4.Cells.Item($r, 5) $donorQ4 $row[3]
    Set-TextValue $q4.Cells.Item($r, 6) $donorQ4 $row[4]
    Set-TextValue $q4.Cells.Item($r, 7) $donorQ4 $row[5]
    $q4.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ===========================================================================
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q4 at row 2,
#    shifting the existing quarters down by one row.
# ===========================================================================

# Prime the style for what will become the new last row (row 6) by copying
# the format of the current last row (row 5) down one row first.
$summary.Cells.Item(5, 1).Copy()
$summary.Cells.Item(6, 1).PasteSpecial(-4122)

for ($row = 5; $row -ge 2; $row--) {
    $bText = $summary.Cells.Item($row, 2).Text
    $cText = $summary.Cells.Item($row, 3).Text
    $dText = $summary.Cells.Item($row, 4).Text
    $summary.Cells.Item($row + 1, 1).Value = ($row - 1)
    $summary.Cells.Item($row + 1, 2).Value = $bText
    $summary.Cells.Item($row + 1, 3).Value = $cText
    $summary.Cells.Item($row + 1, 4).Value = $dText
}

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 5
$summary.Cells.Item(2, 4).Value = 3.65
